$wb = $excel.ActiveWorkbook

# --- Trip sheet: TNC mode balancing-trip fixes ---------------------------
$wsTrip = $wb.Worksheets.Item("Trip")

# New "dorp" categories for multi-passenger TNC rides (rows 6-8, cols M/N)
$wsTrip.Range("M6").Value = 11
$wsTrip.Range("N6").Value = "Driver+Passenger"

$wsTrip.Range("M7").Value = 12
$wsTrip.Range("N7").Value = "Driver+2Passengers"

$wsTrip.Range("M8").Value = 13
$wsTrip.Range("N8").Value = "Driver+3+Passengers"

# Balancing-trip fix: mode code on row 9 changes from "Other" to "PRS"
$wsTrip.Range("J9").Value = "PRS"

# --- Restore the view state: user ends up on the Trip tab ----------------
$wsTour = $wb.Worksheets.Item("Tour")
$wsTour.Activate()
$wsTour.Range("H10").Select()

$wsTrip.Activate()
$wsTrip.Range("J9").Select()
